$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits to the "Stats Lab" module sequence ---
# E6 was "Intro to R part 4 cont..." -> now shows the Stats Lab 5 entry
$ws.Range("E6").Value = "Stats Lab 5: Model validation and presentation"

# E10 was "Stats Lab 2: Generalized linear modelling I" -> becomes an
# "Independent Study" slot (same look as the other Independent Study cells,
# i.e. left aligned with no shaded fill) - copy formatting from E12, which
# already carries that look
$ws.Range("E10").Value = "Independent Study"
$ws.Range("E12").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats

# D12 was "Stats Lab 3: Modelling Your Data I" -> becomes Stats Lab 2
$ws.Range("D12").Value = "Stats Lab 2: Generalized linear modelling I"

# D13 was "Stats Lab 4: Modelling Your Data II" -> renamed to Stats Lab 3
$ws.Range("D13").Value = "Stats Lab 3: Modelling Your Data II"

# D15 stays "Stats Lab 5: Model validation and presentation" (text unchanged)
$ws.Range("D15").Value = "Stats Lab 5: Model validation and presentation"

# F12 stays "Stats Assignment 6 DUE" (text unchanged)
$ws.Range("F12").Value = "Stats Assignment 6 DUE"

# --- Selection change recorded in the sheet view ---
$ws.Range("C6").Select()
